$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.875.89'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '2.440.92'
$ws.Range('E3').Value = '  -1.40%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '558.73'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '162.07'
$ws.Range('E6').Value = '  -1.58%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.511'
$ws.Range('E8').Value = '  -0.62%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.170'
$ws.Range('E9').Value = '  +7.49%  '
$ws.Range('E10').Value = '  -2.22%  '
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.61'
$ws.Range('E12').Value = '  -5.25%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000178'
$ws.Range('E13').Value = '  +3.56%  '
$ws.Range('D14').Value = '68.761.53'
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').Value = '2.888.23'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '23.26'
$ws.Range('E16').Value = '  -1.80%  '
$ws.Range('D17').Value = '2.437.83'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.60'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '339.55'
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('E20').Value = '  -0.31%  '
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.94'
$ws.Range('E22').Value = '  +2.10%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '67.12'
$ws.Range('E24').Value = '  +0.32%  '
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('D26').Value = '2.566.78'
$ws.Range('E26').Value = '  -1.22%  '
$ws.Range('E27').Value = '  -0.82%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.968'
$ws.Range('E28').Value = '  -2.72%  '
$ws.Range('D29').Value = '0.0₃0821'
$ws.Range('E29').Value = '  -0.92%  '
$ws.Range('E30').Value = '  -1.59%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  +1.14%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '430.06'
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('E34').Value = '  -2.20%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '159.83'
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('E38').Value = '  +0.55%  '
$ws.Range('E39').Value = '  -2.59%  '
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('E41').Value = '  +1.86%  '
$ws.Range('E42').Value = '  -2.76%  '
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('E44').Value = '  -1.05%  '
$ws.Range('E45').Value = '  -2.19%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '130.00'
$ws.Range('E46').Value = '  -0.73%  '
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('E49').Value = '  -1.25%  '
$ws.Range('B50').Value = 'BitgetToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.15'
$ws.Range('E50').Value = '  +2.79%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0922'
$ws.Range('E51').Value = '  +0.19%  '
